$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "A" column lag regression coefficients
$ws.Range("B2").Value = "-0.372***"
$ws.Range("B3").Value = "-3.464***"

# Update the "C" column lag regression coefficients
$ws.Range("C2").Value = "0.01*"
$ws.Range("C3").Value = "-0.808***"
